$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.961.42"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "2.388.76"
$ws.Range("E3").Value = "  +2.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").Value = "2.388.58"
$ws.Range("E9").Value = "  +2.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.37%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.73%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.816.26"
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "61.921.30"
$ws.Range("E16").Value = "  +3.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.58%  "

$ws.Range("D18").Value = "2.393.19"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  +6.28%  "

$ws.Range("E28").Value = "  +9.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.31%  "

$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("D31").Value = "0.0₃0765"
$ws.Range("E31").Value = "  +5.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.391"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "341.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0963"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0517"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.28%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.583"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.67%  "

$ws.Range("D51").Value = "0.0₆0221"
$ws.Range("E51").Value = "  -2.34%  "
